$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume snapshot (and two row swaps for
# Bittensor/InternetComputer(DFINITY) and EthereumClassic/RenderToken).
# Numeric-looking Price values are prefixed with a literal apostrophe so
# Excel stores them as text (preserving exact formatting) instead of
# silently converting them to floating point numbers.

$ws.Range("D2").Value = "88.727.51"

$ws.Range("D3").Value = "3.378.33"
$ws.Range("E3").Value = "  +5.91%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'222.87"
$ws.Range("E5").Value = "  +5.99%  "

$ws.Range("D6").Value = "'663.94"
$ws.Range("E6").Value = "  +5.67%  "

$ws.Range("D7").Value = "'0.349"
$ws.Range("E7").Value = "  +25.35%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "'0.620"
$ws.Range("E9").Value = "  +5.35%  "

$ws.Range("D10").Value = "3.378.83"
$ws.Range("E10").Value = "  +6.00%  "

$ws.Range("D11").Value = "'0.616"
$ws.Range("E11").Value = "  +4.46%  "

$ws.Range("D12").Value = "'0.0000278"
$ws.Range("E12").Value = "  +6.43%  "

$ws.Range("E13").Value = "  +2.36%  "

$ws.Range("D14").Value = "'35.62"
$ws.Range("E14").Value = "  +11.32%  "

$ws.Range("D15").Value = "3.988.86"
$ws.Range("E15").Value = "  +5.65%  "

$ws.Range("D16").Value = "'5.51"
$ws.Range("E16").Value = "  +4.31%  "

$ws.Range("D17").Value = "88.308.56"
$ws.Range("E17").Value = "  +9.82%  "

$ws.Range("D18").Value = "3.365.71"
$ws.Range("E18").Value = "  +5.46%  "

$ws.Range("D19").Value = "'14.90"
$ws.Range("E19").Value = "  +4.28%  "

$ws.Range("D20").Value = "'3.25"
$ws.Range("E20").Value = "  +7.35%  "

$ws.Range("D21").Value = "'473.36"
$ws.Range("E21").Value = "  +6.97%  "

$ws.Range("D22").Value = "'5.76"
$ws.Range("E22").Value = "  +9.93%  "

$ws.Range("D23").Value = "'9.35"
$ws.Range("E23").Value = "  +1.59%  "

$ws.Range("D24").Value = "'13.67"
$ws.Range("E24").Value = "  +25.36%  "

$ws.Range("D25").Value = "'7.57"
$ws.Range("E25").Value = "  +9.07%  "

$ws.Range("D26").Value = "'5.57"
$ws.Range("E26").Value = "  +18.27%  "

$ws.Range("E27").Value = "  +4.77%  "

$ws.Range("D28").Value = "'79.99"
$ws.Range("E28").Value = "  +4.70%  "

$ws.Range("D29").Value = "'0.206"
$ws.Range("E29").Value = "  +65.88%  "

$ws.Range("D30").Value = "'0.0000131"
$ws.Range("E30").Value = "  +6.38%  "

$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'609.54"
$ws.Range("E32").Value = "  +8.77%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'9.47"
$ws.Range("E33").Value = "  +5.54%  "

$ws.Range("D34").Value = "'1.58"
$ws.Range("E34").Value = "  +9.30%  "

$ws.Range("D35").Value = "'0.991"
$ws.Range("E35").Value = "  -0.74%  "

$ws.Range("D36").Value = "'2.10"
$ws.Range("E36").Value = "  +5.35%  "

$ws.Range("D37").Value = "'0.153"
$ws.Range("E37").Value = "  +0.55%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'24.34"
$ws.Range("E38").Value = "  +5.62%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'6.97"
$ws.Range("E39").Value = "  +23.52%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.20"
$ws.Range("E40").Value = "  +21.08%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.428"
$ws.Range("E41").Value = "  +5.25%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").Value = "'21.69"
$ws.Range("E43").Value = "  +4.33%  "

$ws.Range("D44").Value = "'3.14"
$ws.Range("E44").Value = "  +14.97%  "

$ws.Range("D45").Value = "'194.74"
$ws.Range("E45").Value = "  +2.31%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'158.14"
$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("D48").Value = "'48.08"
$ws.Range("E48").Value = "  +12.24%  "

$ws.Range("D49").Value = "'1.41"
$ws.Range("E49").Value = "  +8.41%  "

$ws.Range("D50").Value = "'0.807"
$ws.Range("E50").Value = "  +2.92%  "

$ws.Range("D51").Value = "'26.99"
$ws.Range("E51").Value = "  +5.72%  "
